$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 23: 2026-01-16 Kaspa buy.
# The date column stores plain text like "MM/DD/YYYY" (matching the rest of
# the sheet), so force text entry with a leading apostrophe to stop Excel's
# autodetect from turning it into a date serial, then restore the default
# "Normal" style so no extra per-cell formatting is left behind.
$ws.Range("A23").Value = "'01/16/2026"
$ws.Range("A23").Style = "Normal"

$ws.Range("B23").Value = 1076.170999999998
$ws.Range("C23").Value = 0.04599640763410282
$ws.Range("D23").Value = 50
